# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" stats sheet: new data pull for several
# countries, which also reshuffles a few adjacent rows because the sheet
# is kept sorted by total cases (column B) descending. The row swaps
# (España/Colombia, Turquia/Filipinas, Curazao/Islas Feroe,
# Nueva Caledonia/Santa Lucia) are expressed here as direct per-cell
# writes of the new resulting values at each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 19:41"
$ws.Range("B4").Value = 7692391
$ws.Range("C4").Value = 12747
$ws.Range("D4").Value = 4910328
$ws.Range("E4").Value = 2566785
$ws.Range("G4").Value = 246
$ws.Range("H4").Value = 215278
$ws.Range("B5").Value = 6737544
$ws.Range("C5").Value = 55471
$ws.Range("D5").Value = 5718802
$ws.Range("E5").Value = 914563
$ws.Range("G5").Value = 579
$ws.Range("H5").Value = 104179
$ws.Range("A8").Value = "España"
$ws.Range("B8").Value = 865631
$ws.Range("C8").Value = 12793
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("G8").Value = 261
$ws.Range("H8").Value = 32486
$ws.Range("A9").Value = "Colombia"
$ws.Range("B9").Value = 862158
$ws.Range("D9").Value = 766300
$ws.Range("E9").Value = 69014
$ws.Range("H9").Value = 26844
$ws.Range("A22").Value = "Turquia"
$ws.Range("B22").Value = 327557
$ws.Range("C22").Value = 1511
$ws.Range("D22").Value = 287599
$ws.Range("E22").Value = 31405
$ws.Range("G22").Value = 55
$ws.Range("H22").Value = 8553
$ws.Range("A23").Value = "Filipinas"
$ws.Range("B23").Value = 326833
$ws.Range("C23").Value = 2093
$ws.Range("D23").Value = 273313
$ws.Range("E23").Value = 47655
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = 5865
$ws.Range("B26").Value = 305869
$ws.Range("C26").Value = 1212
$ws.Range("E26").Value = 32542
$ws.Range("G26").Value = 11
$ws.Range("H26").Value = 9627
$ws.Range("B27").Value = 276439
$ws.Range("C27").Value = 4130
$ws.Range("D27").Value = 210845
$ws.Range("E27").Value = 63810
$ws.Range("G27").Value = 27
$ws.Range("H27").Value = 1784
$ws.Range("B74").Value = 38973
$ws.Range("C74").Value = 424
$ws.Range("E74").Value = 13798
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 1811
$ws.Range("B137").Value = 4252
$ws.Range("C137").Value = 739
$ws.Range("E137").Value = 973
$ws.Range("A179").Value = "Curazao"
$ws.Range("B179").Value = 476
$ws.Range("C179").Value = 14
$ws.Range("D179").Value = 262
$ws.Range("E179").Value = 213
$ws.Range("H179").Value = 1
$ws.Range("A180").Value = "Islas Feroe"
$ws.Range("B180").Value = 475
$ws.Range("D180").Value = 452
$ws.Range("E180").Value = 23
$ws.Range("H180").Value = 0
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"
